$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.099.24"
$ws.Range("E2").Value = "  -2.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.239.77"
$ws.Range("E3").Value = "  -2.57%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.46"
$ws.Range("E5").Value = "  -3.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.61"
$ws.Range("E6").Value = "  -8.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.578"
$ws.Range("E7").Value = "  -4.85%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.234.66"
$ws.Range("E9").Value = "  -2.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.589"
$ws.Range("E10").Value = "  -5.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "51.51"
$ws.Range("E11").Value = "  -13.59%  "
$ws.Range("E12").Value = "  -3.44%  "
$ws.Range("E13").Value = "  -5.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.74"
$ws.Range("E14").Value = "  -4.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.759.85"
$ws.Range("E15").Value = "  -2.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.239.70"
$ws.Range("E16").Value = "  -2.72%  "
$ws.Range("E17").Value = "  -2.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.876.40"
$ws.Range("E18").Value = "  -2.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.01"
$ws.Range("E19").Value = "  -4.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.92"
$ws.Range("E20").Value = "  -2.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.934"
$ws.Range("E21").Value = "  -3.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "372.31"
$ws.Range("E22").Value = "  -1.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.08"
$ws.Range("E23").Value = "  +4.67%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.19"
$ws.Range("E24").Value = "  -1.41%  "
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.96"
$ws.Range("E25").Value = "  -3.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.59"
$ws.Range("E26").Value = "  -6.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.12"
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.63"
$ws.Range("E28").Value = "  -2.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.96"
$ws.Range("E30").Value = "  -6.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.93"
$ws.Range("E31").Value = "  -4.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "615.21"
$ws.Range("E32").Value = "  -4.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.46"
$ws.Range("E33").Value = "  -5.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.97"
$ws.Range("E34").Value = "  -3.45%  "
$ws.Range("E35").Value = "  -3.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "56.34"
$ws.Range("E36").Value = "  -6.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "34.90"
$ws.Range("E38").Value = "  -5.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.368"
$ws.Range("E39").Value = "  -7.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0703"
$ws.Range("E41").Value = "  -3.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.58"
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("E43").Value = "  -4.80%  "
$ws.Range("E44").Value = "  +4.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.816.18"
$ws.Range("E45").Value = "  -3.59%  "
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("E47").Value = "  -4.49%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.00"
$ws.Range("E48").Value = "  +0.51%  "
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.52"
$ws.Range("E49").Value = "  -7.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "136.34"
$ws.Range("E50").Value = "  +1.56%  "
$ws.Range("E51").Value = "  -3.69%  "
